$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add two new localization rows (54 and 55) to the "General" sheet, matching
# the rows already present for the other cwl_* log/warn id+filter+text_JP+
# text entries (rows 50-53). Column B ("filter") is always left blank.
# ---------------------------------------------------------------------------

# --- Row 54 -----------------------------------------------------------------
# Copy formatting (cell styles) from the row-50 quadruplet, which uses the
# exact same style set we need here (s="4" / s="9" / s="10" / s="10").
$ws.Range("A50").Copy()
$ws.Range("A54").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B50").Copy()
$ws.Range("B54").PasteSpecial(-4122)
$ws.Range("C50").Copy()
$ws.Range("C54").PasteSpecial(-4122)
$ws.Range("D50").Copy()
$ws.Range("D54").PasteSpecial(-4122)

$ws.Range("A54").Value = "cwl_log_custom_trait"
$ws.Range("C54").Value = "qualified custom trait id: {0}, type: {1}"

$ws.Range("D54").Value = "已加载自定义特征: {0}, 限定类型: {1}"
$r = $ws.Range("D54").Characters(1, 8)
$r.Font.Name = "宋体"
$r.Font.Size = 15.8
$r = $ws.Range("D54").Characters(9, 7)
$r.Font.Name = "Cascadia Code"
$r.Font.Size = 15.8
$r = $ws.Range("D54").Characters(16, 4)
$r.Font.Name = "宋体"
$r.Font.Size = 15.8
$r = $ws.Range("D54").Characters(20, 5)
$r.Font.Name = "Cascadia Code"
$r.Font.Size = 15.8

# --- Row 55 -----------------------------------------------------------------
# Row 51 carries the matching style set for row 55 (s="2" / s="9" / s="10" /
# s="10").
$ws.Range("A51").Copy()
$ws.Range("A55").PasteSpecial(-4122)
$ws.Range("B51").Copy()
$ws.Range("B55").PasteSpecial(-4122)
$ws.Range("C51").Copy()
$ws.Range("C55").PasteSpecial(-4122)
$ws.Range("D51").Copy()
$ws.Range("D55").PasteSpecial(-4122)

$ws.Range("A55").Value = "cwl_warn_qualify_trait"
$ws.Range("C55").Value = "failed to qualify custom trait id: {0}, on card: {1}"

$ws.Range("D55").Value = "无法限定自定义特征: {0}, 卡片: {1}"
$r = $ws.Range("D55").Characters(1, 9)
$r.Font.Name = "微软雅黑"
$r.Font.Size = 15.8
$r = $ws.Range("D55").Characters(10, 7)
$r.Font.Name = "Cascadia Code"
$r.Font.Size = 15.8
$r = $ws.Range("D55").Characters(17, 2)
$r.Font.Name = "宋体"
$r.Font.Size = 15.8
$r = $ws.Range("D55").Characters(19, 5)
$r.Font.Name = "Cascadia Code"
$r.Font.Size = 15.8

# Leave the cursor/selection near the new rows, like the authored edit did.
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("D57").Select()
